# Generate Report for Handoff
# Update status from "In Translation" to "Ready for handoff" and bump the
# related handoff timestamps on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-03-22 08:33:42"

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-22 08:33:38"

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-22 08:33:42"
